$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 9 for the CHE1/region breakdown (shifts BRA1 + NYC rows down by one)
$ws1.Rows("9:9").Insert()

$ws1.Range("A9").Value = "CHE1"
$ws1.Range("B9").Value = "region"
$ws1.Range("C9").Value = "data/derived/CHE/CHE_region.RDS"
$ws1.Range("D9").Value = "stratified"
$ws1.Range("E9").Value = "aggregate"

# Append new row 13 for GBR2 region
$ws1.Range("A13").Value = "GBR2"
$ws1.Range("B13").Value = "region"
$ws1.Range("C13").Value = "data/derived/UK/GBR_regions.RDS"
$ws1.Range("D13").Value = "marginal"
$ws1.Range("E13").Value = "linelist"

# Rename ESP1 -> ESP1-2 (rows 2-3, column A)
$ws1.Range("A2").Value = "ESP1-2"
$ws1.Range("A3").Value = "ESP1-2"

# Append new row 14 for GBR2 ageband
$ws1.Range("A14").Value = "GBR2"
$ws1.Range("B14").Value = "ageband"
$ws1.Range("C14").Value = "data/derived/UK/GBR2_agebands.RDS"
$ws1.Range("D14").Value = "marginal"
$ws1.Range("E14").Value = "linelist"

# Update selection to match final state (C14 selected)
[void]$ws1.Range("C14").Select()
